# Updated PCMC code of payments and reports
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - login details now point to PCMC test site / user "sagar.m"
$ws.Range("A2").Value = "sagar.m"
$ws.Range("B2").Value = 123
$ws.Range("C2").Value = "http://testpcmc.ptaxcollection.com:8080/Pages/Login.aspx"

# Row 10 - BMC -> TRG sector, counts updated, node code updated
$ws.Range("C10").Value = "131"
$ws.Range("A10").Value = "TRG"
$ws.Range("B10").Value = 10

# Row 11
$ws.Range("A11").Value = "TRG"
$ws.Range("B11").Value = 6
$ws.Range("C11").Value = "92"

# Row 12 - only the node code changes
$ws.Range("C12").Value = "20"

# Row 13
$ws.Range("A13").Value = "TRG"
$ws.Range("B13").Value = 6
$ws.Range("C13").Value = "93"

# Update the active cell/selection to A14
$ws.Range("A14").Select()
